$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp text (A1)
$ws.Range("A1").Value = "Datos actualizados a 26 de Mayo de 2020 a las 15:05"

# Helper data: row number -> @(Country, Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
$rows = @{
    4   = @("Estados Unidos", 1708079, 1853, 464728, 1143510, 0, 36, 99841)
    11  = @("Alemania", 180830, 41, 162000, 10397, 0, 5, 8433)
    13  = @("India", 146498, 1548, 61151, 81160, 0, 15, 4187)
    28  = @("Suecia", 34440, 597, 4971, 25344, 0, 96, 4125)
    31  = @("Portugal", 31007, 219, 18096, 11569, 0, 12, 1342)
    50  = @("Serbia", 11227, 34, 6067, 4921, 0, 0, 239)
    51  = @("Corea del Sur", 11225, 19, 10275, 681, 0, 2, 269)
    65  = @("Ghana", 6964, 156, 2097, 4835, 0, 0, 32)
    78  = @("Tayikistan", 3266, 166, 1417, 1802, 0, 1, 47)
    79  = @("Senegal", 3161, 31, 1565, 1560, 0, 0, 36)
    103 = @("Sri Lanka", 1206, 24, 712, 484, 0, 0, 10)
    207 = @("Groenlandia", 12, 0, 11, 1, 0, 0, 0)
    208 = @("Islas Turcas y Caicos", 12, 0, 10, 1, 0, 0, 1)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = $vals[3]
    $ws.Cells.Item($r, 5).Value = $vals[4]
    $ws.Cells.Item($r, 6).Value = $vals[5]
    $ws.Cells.Item($r, 7).Value = $vals[6]
    $ws.Cells.Item($r, 8).Value = $vals[7]
}
